# Auto-generated script applying the Zalera_Profits.xlsx data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 265.58334
$ws.Cells.Item(5, 9).Value = 231.42857
$ws.Cells.Item(5, 11).Value = 231.42857
$ws.Cells.Item(5, 13).Value = -116.42857
$ws.Cells.Item(40, 8).Value = 3086.3333
$ws.Cells.Item(40, 9).Value = 2710.0908
$ws.Cells.Item(40, 10).Value = 3677.5715
$ws.Cells.Item(40, 11).Value = 2710.0908
$ws.Cells.Item(40, 12).Value = 3677.5715
$ws.Cells.Item(40, 13).Value = -2535.0908
$ws.Cells.Item(40, 14).Value = -4027.5715
$ws.Cells.Item(135, 8).Value = 5212.077
$ws.Cells.Item(135, 9).Value = 4026.1667
$ws.Cells.Item(135, 11).Value = 36235.5003
$ws.Cells.Item(135, 13).Value = -33700.5003
$ws.Cells.Item(137, 8).Value = 5065
$ws.Cells.Item(137, 9).Value = 1461.8422
$ws.Cells.Item(137, 10).Value = 14845
$ws.Cells.Item(137, 11).Value = 4385.5266
$ws.Cells.Item(137, 12).Value = 44535
$ws.Cells.Item(137, 13).Value = -1835.5266
$ws.Cells.Item(137, 14).Value = -49635
$ws.Cells.Item(138, 8).Value = 3847.5483
$ws.Cells.Item(138, 9).Value = 1535
$ws.Cells.Item(138, 10).Value = 4793.591
$ws.Cells.Item(138, 11).Value = 4605
$ws.Cells.Item(138, 12).Value = 14380.773
$ws.Cells.Item(138, 13).Value = 535
$ws.Cells.Item(138, 14).Value = -24660.773
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 7493.375
$ws.Cells.Item(45, 9).Value = 8607.385
$ws.Cells.Item(45, 11).Value = 8607.385
$ws.Cells.Item(45, 13).Value = -8230.385
$ws.Cells.Item(61, 8).Value = 4356.0234
$ws.Cells.Item(61, 9).Value = 4149.6665
$ws.Cells.Item(61, 11).Value = 4149.6665
$ws.Cells.Item(61, 13).Value = -3937.6665
$ws.Cells.Item(88, 8).Value = 1713.4546
$ws.Cells.Item(88, 9).Value = 1691.8572
$ws.Cells.Item(88, 10).Value = 1751.25
$ws.Cells.Item(88, 11).Value = 1691.8572
$ws.Cells.Item(88, 12).Value = 1751.25
$ws.Cells.Item(88, 13).Value = -1285.8572
$ws.Cells.Item(88, 14).Value = -2563.25
$ws.Cells.Item(91, 8).Value = 1713.4546
$ws.Cells.Item(91, 9).Value = 1691.8572
$ws.Cells.Item(91, 10).Value = 1751.25
$ws.Cells.Item(91, 11).Value = 1691.8572
$ws.Cells.Item(91, 12).Value = 1751.25
$ws.Cells.Item(91, 13).Value = -287.8571999999999
$ws.Cells.Item(91, 14).Value = -4559.25
$ws.Cells.Item(132, 8).Value = 3448.68
$ws.Cells.Item(132, 9).Value = 2889.9487
$ws.Cells.Item(132, 10).Value = 5429.636
$ws.Cells.Item(132, 11).Value = 8669.846099999999
$ws.Cells.Item(132, 12).Value = 16288.908
$ws.Cells.Item(132, 13).Value = -6139.846099999999
$ws.Cells.Item(132, 14).Value = -21348.908
$ws.Cells.Item(136, 8).Value = 4356.0234
$ws.Cells.Item(136, 9).Value = 4149.6665
$ws.Cells.Item(136, 11).Value = 12448.9995
$ws.Cells.Item(136, 13).Value = -9898.999500000002
$ws.Cells.Item(137, 8).Value = 120000
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(141, 8).Value = 114949.164
$ws.Cells.Item(141, 10).Value = 114949.164
$ws.Cells.Item(141, 12).Value = 114949.164
$ws.Cells.Item(141, 14).Value = -125309.164
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 79735.39999999999
$ws.Cells.Item(122, 9).Value = 168357.28
$ws.Cells.Item(122, 11).Value = 505071.84
$ws.Cells.Item(122, 13).Value = -502621.84
$ws.Cells.Item(132, 8).Value = 82558.484
$ws.Cells.Item(132, 9).Value = 57005.22
$ws.Cells.Item(132, 10).Value = 246099.4
$ws.Cells.Item(132, 11).Value = 171015.66
$ws.Cells.Item(132, 12).Value = 738298.2
$ws.Cells.Item(132, 13).Value = -168485.66
$ws.Cells.Item(132, 14).Value = -743358.2
$ws.Cells.Item(134, 8).Value = 3414.4187
$ws.Cells.Item(134, 10).Value = 6458.846
$ws.Cells.Item(134, 12).Value = 19376.538
$ws.Cells.Item(134, 14).Value = -24446.538
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 92253.5
$ws.Cells.Item(37, 10).Value = 92253.5
$ws.Cells.Item(37, 12).Value = 276760.5
$ws.Cells.Item(37, 14).Value = -276984.5
$ws.Cells.Item(48, 8).Value = 3000
$ws.Cells.Item(48, 10).Value = 3000
$ws.Cells.Item(48, 12).Value = 9000
$ws.Cells.Item(48, 14).Value = -9500
$ws.Cells.Item(61, 8).Value = 242.2973
$ws.Cells.Item(61, 10).Value = 285.13333
$ws.Cells.Item(61, 12).Value = 855.39999
$ws.Cells.Item(61, 14).Value = -1285.39999
$ws.Cells.Item(62, 8).Value = 1954.9788
$ws.Cells.Item(65, 8).Value = 1954.9788
$ws.Cells.Item(113, 8).Value = 1456.7241
$ws.Cells.Item(113, 10).Value = 1610.8889
$ws.Cells.Item(113, 12).Value = 4832.6667
$ws.Cells.Item(113, 14).Value = -9172.6667
$ws.Cells.Item(131, 8).Value = 8777375
$ws.Cells.Item(131, 10).Value = 6053.697
$ws.Cells.Item(131, 12).Value = 18161.091
$ws.Cells.Item(131, 14).Value = -28241.091
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3145
$ws.Cells.Item(80, 9).Value = 2872.5
$ws.Cells.Item(80, 10).Value = 3387.2222
$ws.Cells.Item(80, 11).Value = 2872.5
$ws.Cells.Item(80, 12).Value = 3387.2222
$ws.Cells.Item(80, 13).Value = -1874.5
$ws.Cells.Item(80, 14).Value = -5383.2222
$ws.Cells.Item(83, 8).Value = 3145
$ws.Cells.Item(83, 9).Value = 2872.5
$ws.Cells.Item(83, 10).Value = 3387.2222
$ws.Cells.Item(83, 11).Value = 14362.5
$ws.Cells.Item(83, 12).Value = 16936.111
$ws.Cells.Item(83, 13).Value = -9370.5
$ws.Cells.Item(83, 14).Value = -26920.111
$ws.Cells.Item(113, 8).Value = 38705.5
$ws.Cells.Item(113, 9).Value = 38705.5
$ws.Cells.Item(113, 11).Value = 38705.5
$ws.Cells.Item(113, 13).Value = -36535.5
$ws.Cells.Item(138, 8).Value = 110429
$ws.Cells.Item(138, 10).Value = 110429
$ws.Cells.Item(138, 12).Value = 110429
$ws.Cells.Item(138, 14).Value = -120709
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3828.2
$ws.Cells.Item(68, 9).Value = 3098.4
$ws.Cells.Item(68, 11).Value = 3098.4
$ws.Cells.Item(68, 13).Value = -2349.4
$ws.Cells.Item(71, 8).Value = 3828.2
$ws.Cells.Item(71, 9).Value = 3098.4
$ws.Cells.Item(71, 11).Value = 15492
$ws.Cells.Item(71, 13).Value = -11748
$ws.Cells.Item(136, 8).Value = 2890.9285
$ws.Cells.Item(136, 9).Value = 2299.8518
$ws.Cells.Item(136, 11).Value = 6899.555399999999
$ws.Cells.Item(136, 13).Value = -4349.555399999999
$ws.Cells.Item(137, 8).Value = 128888.664
$ws.Cells.Item(137, 10).Value = 128888.664
$ws.Cells.Item(137, 12).Value = 128888.664
$ws.Cells.Item(137, 14).Value = -139088.664
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 10).Value = 1000
$ws.Cells.Item(9, 12).Value = 1000
$ws.Cells.Item(9, 14).Value = -1280
$ws.Cells.Item(132, 9).Value = 2619.3057
$ws.Cells.Item(132, 11).Value = 7857.9171
$ws.Cells.Item(132, 13).Value = -5327.9171
$ws.Cells.Item(136, 8).Value = 2537.8333
$ws.Cells.Item(136, 9).Value = 1359.4546
$ws.Cells.Item(136, 11).Value = 4078.3638
$ws.Cells.Item(136, 13).Value = -1528.3638
